$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.667213916778564
$ws.Range("B1").Value = 1.916872262954712
$ws.Range("C1").Value = 1.947595357894897
$ws.Range("D1").Value = 2.494866371154785
$ws.Range("E1").Value = 3.612186431884766
